# Apply "Update countries & provincias Spain" data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 1; Values = @("Datos actualizados a 24 de Agosto de 2020 a las 01:55") },
    @{ Row = 4; Values = @("Estados Unidos", 5873237, 31809, 3162810, 2529847, 0, 406, 180580) },
    @{ Row = 5; Values = @("Brasil", 3605783, 23085, 2709638, 781373, 0, 495, 114772) },
    @{ Row = 15; Values = @("Argentina", 342154, 5352, 251400, 83769, 0, 137, 6985) },
    @{ Row = 22; Values = @("Francia", 242899, 4897, 84973, 127413, 0, 1, 30513) },
    @{ Row = 23; Values = @("Alemania", 234489, 632, 209600, 15557, 0, 1, 9332) },
    @{ Row = 36; Values = @("Panama", 86900, 1420, 61420, 23588, 0, 14, 1892) },
    @{ Row = 37; Values = @("Suecia", 86068, 0, 0, 0, 0, 0, 5810) },
    @{ Row = 74; Values = @("Chequia", 21923, 133, 16125, 5386, 0, 1, 412) },
    @{ Row = 102; Values = @("Luxemburgo", 7775, 13, 6969, 682, 0, 0, 124) },
    @{ Row = 104; Values = @("Maldivas", 6779, 119, 4222, 2531, 0, 0, 26) },
    @{ Row = 105; Values = @("Zimbabue", 5930, 37, 4872, 903, 0, 2, 155) },
    @{ Row = 107; Values = @("Malaui", 5414, 32, 3012, 2234, 0, 0, 168) },
    @{ Row = 113; Values = @("Montenegro", 4343, 30, 3356, 903, 0, 0, 84) },
    @{ Row = 115; Values = @("Suazilandia", 4225, 36, 2898, 1242, 0, 2, 85) },
    @{ Row = 117; Values = @("Cuba", 3682, 65, 3044, 547, 0, 2, 91) },
    @{ Row = 118; Values = @("Surinam", 3607, 38, 2688, 861, 0, 1, 58) },
    @{ Row = 140; Values = @("Sierra Leona", 1992, 12, 1550, 373, 0, 0, 69) },
    @{ Row = 147; Values = @("Uruguay", 1527, 6, 1276, 209, 0, 0, 42) },
    @{ Row = 163; Values = @("Guyana", 955, 30, 490, 434, 0, 0, 31) },
    @{ Row = 168; Values = @("Belice", 686, 18, 44, 636, 0, 0, 6) },
    @{ Row = 174; Values = @("Papua Nueva Guinea", 401, 40, 232, 165, 0, 0, 4) },
    @{ Row = 175; Values = @("Islas Feroe", 384, 0, 300, 84, 0, 0, 0) },
    @{ Row = 176; Values = @("Islas Turcas y Caicos", 383, 36, 102, 279, 0, 0, 2) },
    @{ Row = 177; Values = @("San Martin (Parte Holandesa)", 368, 0, 147, 204, 0, 0, 17) },
    @{ Row = 202; Values = @("Timor Oriental", 26, 0, 25, 1, 0, 0, 0) },
    @{ Row = 203; Values = @("Santa Lucia", 26, 0, 25, 1, 0, 0, 0) }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $vals = $r.Values
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c).Value = $vals[$c - 1]
    }
}
